$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table with the latest scraped
# values. Several coins (OKB/Hedera, Dai/Fetch.AI/TheGraph, Maker/ThetaToken)
# also changed rank order, so their Name/Link/Price/Volume cells are updated
# in place to reflect the new row positions.
#
# Price cells (column D) are assigned with a leading apostrophe so that
# numeric-looking text (e.g. "573.12", "0.185") is stored as literal text,
# matching the source data, instead of being auto-converted to a number.
$ws.Range('D2').Value = "'69.120.20"
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = "'3.521.29"
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'573.12"
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = "'183.78"
$ws.Range('E6').Value = '  -2.50%  '
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('D8').Value = "'3.515.41"
$ws.Range('E8').Value = '  -1.74%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = "'0.185"
$ws.Range('E10').Value = '  +3.95%  '
$ws.Range('D11').Value = "'0.641"
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').Value = "'53.96"
$ws.Range('E12').Value = '  -3.68%  '
$ws.Range('D13').Value = "'0.0000303"
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').Value = "'4.090.09"
$ws.Range('E15').Value = '  -1.58%  '
$ws.Range('D16').Value = "'19.36"
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').Value = "'3.520.54"
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = "'69.097.96"
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').Value = "'12.54"
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('E20').Value = '  -1.09%  '
$ws.Range('D21').Value = "'539.13"
$ws.Range('E21').Value = '  +13.73%  '
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').Value = "'20.70"
$ws.Range('E23').Value = '  +8.70%  '
$ws.Range('D24').Value = "'5.01"
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').Value = "'94.75"
$ws.Range('E26').Value = '  +6.77%  '
$ws.Range('D27').Value = "'10.97"
$ws.Range('E27').Value = '  -0.87%  '
$ws.Range('D28').Value = "'2.93"
$ws.Range('E28').Value = '  -3.87%  '
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('D30').Value = "'31.57"
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('E31').Value = '  -5.69%  '
$ws.Range('D32').Value = "'12.68"
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = "'64.32"
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'0.114"
$ws.Range('E34').Value = '  -4.48%  '
$ws.Range('D35').Value = "'573.44"
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').Value = "'38.07"
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Value = "'0.401"
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = "'1.00"
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = "'3.06"
$ws.Range('E39').Value = '  +5.00%  '
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('D41').Value = "'0.135"
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('D42').Value = "'3.11"
$ws.Range('E42').Value = '  -3.19%  '
$ws.Range('E43').Value = '  -4.39%  '
$ws.Range('D44').Value = "'3.55"
$ws.Range('E44').Value = '  +6.98%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = "'2.97"
$ws.Range('E45').Value = '  -4.91%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = "'3.207.08"
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = "'0.0440"
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = "'9.16"
$ws.Range('E48').Value = '  -4.12%  '
$ws.Range('E49').Value = '  -1.94%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = "'136.21"
